$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update landmark coordinate values (columns B, C, D) for rows 2-32
$ws.Cells.Item(2, 2).Value = 2.2483823375766208
$ws.Cells.Item(2, 3).Value = 18.011546205610252
$ws.Cells.Item(2, 4).Value = 15.086522852935371
$ws.Cells.Item(3, 2).Value = 2.4107639775420577
$ws.Cells.Item(3, 3).Value = 20.369763188336471
$ws.Cells.Item(3, 4).Value = 14.100090737896041
$ws.Cells.Item(4, 2).Value = 4.6600184659661625
$ws.Cells.Item(4, 3).Value = 16.941927990659938
$ws.Cells.Item(4, 4).Value = 8.7180985036830307
$ws.Cells.Item(5, 2).Value = 2.168594022029243
$ws.Cells.Item(5, 3).Value = 17.319803417587377
$ws.Cells.Item(5, 4).Value = 10.661397614943111
$ws.Cells.Item(6, 2).Value = -10.6782445138626
$ws.Cells.Item(6, 3).Value = 16.52524269391866
$ws.Cells.Item(6, 4).Value = 12.544157943275824
$ws.Cells.Item(7, 2).Value = -11.123440988617357
$ws.Cells.Item(7, 3).Value = 16.653526688921886
$ws.Cells.Item(7, 4).Value = 9.1575794368873034
$ws.Cells.Item(8, 2).Value = -12.417738656417121
$ws.Cells.Item(8, 3).Value = 17.993513417248838
$ws.Cells.Item(8, 4).Value = 11.516940843786514
$ws.Cells.Item(9, 2).Value = -12.571327297883769
$ws.Cells.Item(9, 3).Value = 18.755564073126521
$ws.Cells.Item(9, 4).Value = 8.7075450391447298
$ws.Cells.Item(10, 2).Value = 40.315200017356041
$ws.Cells.Item(10, 3).Value = 11.14103554976106
$ws.Cells.Item(10, 4).Value = 8.7217749276246828
$ws.Cells.Item(11, 2).Value = 40.195438623577346
$ws.Cells.Item(11, 3).Value = 13.644098858355024
$ws.Cells.Item(11, 4).Value = 9.0868152619816929
$ws.Cells.Item(12, 2).Value = 38.386142761120276
$ws.Cells.Item(12, 3).Value = 14.947572152086259
$ws.Cells.Item(12, 4).Value = 9.3252087193500657
$ws.Cells.Item(13, 2).Value = -43.711614242018548
$ws.Cells.Item(13, 3).Value = 4.4660043246097985
$ws.Cells.Item(13, 4).Value = 20.688447645834287
$ws.Cells.Item(14, 2).Value = -42.671362445426794
$ws.Cells.Item(14, 3).Value = 11.613376667260885
$ws.Cells.Item(14, 4).Value = 19.977697258530807
$ws.Cells.Item(15, 2).Value = -39.831366223563968
$ws.Cells.Item(15, 3).Value = 18.454591628516582
$ws.Cells.Item(15, 4).Value = 19.085111823550118
$ws.Cells.Item(16, 2).Value = -45.053217382054697
$ws.Cells.Item(16, 3).Value = 7.82536437087533
$ws.Cells.Item(16, 4).Value = 12.538433704522019
$ws.Cells.Item(17, 2).Value = -42.082459482324161
$ws.Cells.Item(17, 3).Value = 12.579215181039181
$ws.Cells.Item(17, 4).Value = 11.254700802962256
$ws.Cells.Item(18, 2).Value = -39.108122090382764
$ws.Cells.Item(18, 3).Value = 17.206336369660452
$ws.Cells.Item(18, 4).Value = 10.314825999972834
$ws.Cells.Item(19, 2).Value = 41.127390106883986
$ws.Cells.Item(19, 3).Value = 5.0844281756616523
$ws.Cells.Item(19, 4).Value = 11.226907085891133
$ws.Cells.Item(20, 2).Value = 40.726429801347308
$ws.Cells.Item(20, 3).Value = 8.0786209193035887
$ws.Cells.Item(20, 4).Value = 9.9982325507097194
$ws.Cells.Item(21, 2).Value = 40.648067595255206
$ws.Cells.Item(21, 3).Value = 10.859303029494653
$ws.Cells.Item(21, 4).Value = 9.6126989259103031
$ws.Cells.Item(22, 2).Value = 36.977710520601562
$ws.Cells.Item(22, 3).Value = 14.765515469111204
$ws.Cells.Item(22, 4).Value = 5.6275658245338116
$ws.Cells.Item(23, 2).Value = 38.313921786951475
$ws.Cells.Item(23, 3).Value = 18.724344100627842
$ws.Cells.Item(23, 4).Value = 6.7624971952186792
$ws.Cells.Item(24, 2).Value = 38.17880435243746
$ws.Cells.Item(24, 3).Value = 21.511474965780028
$ws.Cells.Item(24, 4).Value = 7.8175632382706386
$ws.Cells.Item(25, 2).Value = -40.302667457108271
$ws.Cells.Item(25, 3).Value = 14.520728471293964
$ws.Cells.Item(25, 4).Value = 7.0154751073601203
$ws.Cells.Item(26, 2).Value = -38.958274026068253
$ws.Cells.Item(26, 3).Value = 16.416454560567544
$ws.Cells.Item(26, 4).Value = 7.1519139611667963
$ws.Cells.Item(27, 2).Value = 28.856087485849759
$ws.Cells.Item(27, 3).Value = 20.164263973099484
$ws.Cells.Item(27, 4).Value = 30.651281284445307
$ws.Cells.Item(28, 2).Value = 17.610305593874497
$ws.Cells.Item(28, 3).Value = 26.129182443438207
$ws.Cells.Item(28, 4).Value = 31.047711813984741
$ws.Cells.Item(29, 2).Value = 9.9077978452714408
$ws.Cells.Item(29, 3).Value = 25.868034504238139
$ws.Cells.Item(29, 4).Value = 28.684693368088059
$ws.Cells.Item(30, 2).Value = -11.405942093791385
$ws.Cells.Item(30, 3).Value = 23.372500680398026
$ws.Cells.Item(30, 4).Value = 13.90607439441805
$ws.Cells.Item(31, 2).Value = -20.981413768657315
$ws.Cells.Item(31, 3).Value = 37.261185910112907
$ws.Cells.Item(31, 4).Value = 16.597289266642932
$ws.Cells.Item(32, 2).Value = -36.237703276710242
$ws.Cells.Item(32, 3).Value = 26.46669746146544
$ws.Cells.Item(32, 4).Value = 14.081747037571965

# Adjust column widths for columns A and B (closest achievable to 9.85546875 / 12.42578125)
$ws.Columns.Item(1).ColumnWidth = 9.0
$ws.Columns.Item(2).ColumnWidth = 11.667
